$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined Results")
$ws.Range("F2").Value = "F1 = 0.84 / 0.85"
$ws.Range("F3").Value = "F1 = 0.84 / 0.85"
$ws.Range("F4").Value = "F1 = 0.84 / 0.85"
$ws.Range("F5").Value = "F1 = 0.84 / 0.85"
$ws.Range("F6").Value = "F1 = 0.82/ 0.21"
$ws.Range("F7").Value = "F1 = 0.82/ 0.21"
$ws.Range("F6").Select()
